$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 85: Hours 1.25 -> 1.75; Notes changes from "Finish 5 small problems" to "Finish 7 small problems"
$ws.Range("C85").Value = 1.75
$ws.Range("D85").Value = "Finish 7 small problems"

# Row 86: add Course (B86) and Milestones (E86)
$ws.Range("B86").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("E86").Value = "begin Lesson 5"

# Update selection to D86 to match the saved cursor position
$ws.Range("D86").Select()
